# BOM.xlsx update — new purchasable parts for NextPCB manufacturing,
# plus assorted corrections (per commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Row 11 (J1 - USB-C receptacle): replace the generic placeholder
# manufacturer/part info with the real sourced part (Bel Fuse / Stewart
# Connector), add a footprint note in PACKAGE and a new NOTES-ish cell
# in F, and expand the description with a sourcing link.
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "    J1"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "USB_C_Receptacle"
$ws.Range("D11").Value = "Footprints:BELFUSE_SS-52400-002"
$ws.Range("E11").Value = "BELFUSE Stewart Connector"
$ws.Range("F11").Value = "SS-52400-003 for standard 1.6 mm PCB thickness"
$ws.Range("G11").Value = "USB 2.0-only Type-C Receptacle connector. See https://www.snapeda.com/parts/SS-52400-002/Bel%20Fuse/view-part/?ref=search&t=BELFUSE_SS-52400-002"

# ---------------------------------------------------------------------
# Row 32 (SW1, SW2 - tactile switches): replace with a real sourced
# part (TE Connectivity USLPT2819 family) and matching links.
# ---------------------------------------------------------------------
$ws.Range("A32").Value = ">  SW1, SW2"
$ws.Range("B32").Value = 2
$ws.Range("C32").Value = "1101NE"
$ws.Range("D32").Value = "USLPT2819 Family"
$ws.Range("E32").Value = "TE Connectivity"
$ws.Range("F32").Value = "USLPT2819 MSLPT (Mini size)"
$ws.Range("G32").Value = "SMD_6x3.5mm_h2.5mm. Please see https://www.snapeda.com/parts/USLPT2819DT2TR/TE%20Connectivity/view-part/?ref=search&t=SMD%20switch"

# ---------------------------------------------------------------------
# Row 34 (U2 - LDO regulator): replace TI LM1117 reference with the
# actual sourced ON Semiconductor part and its datasheet/search link.
# ---------------------------------------------------------------------
$ws.Range("A34").Value = "    U2"
$ws.Range("B34").Value = 1
$ws.Range("C34").Value = "LM1117MPX-3.3_NOPB"
$ws.Range("D34").Value = "LM1117MPX-33NOPB"
$ws.Range("E34").Value = "On Semiconductor"
$ws.Range("F34").Value = "SOT−223 CASE 318H"
$ws.Range("G34").Value = "See https://www.snapeda.com/parts/LM1117MPX-33NOPB/ON%20Semiconductor/view-part/?welcome=home&ref=search&t=LM1117MPX-3.3_NOPB"

# ---------------------------------------------------------------------
# Highlight the three corrected/newly-sourced rows in yellow so they
# stand out as "verified for NextPCB purchasing".
# ---------------------------------------------------------------------
$ws.Range("A11:G11").Interior.Color = 65535
$ws.Range("A32:G32").Interior.Color = 65535
$ws.Range("A34:G34").Interior.Color = 65535

# ---------------------------------------------------------------------
# Column layout: VALUE (B) and PACKAGE->MANUF. PART No area needed more
# room once the new sourcing notes were added to column C, so split the
# previously-shared B:C width and widen C.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 19.5703125

# ---------------------------------------------------------------------
# Selection cursor moved while reviewing the new J1 footprint note.
# ---------------------------------------------------------------------
$ws.Range("F12").Select()
